$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 402; this shifts the existing rows 402-437
# down to 404-439, preserving all of their existing data and formatting.
$ws.Rows.Item(402).Resize(2).Insert()

# Row 402: new "Primera" record dated 2022-08-10 (serial 44783)
$ws.Cells.Item(402, 1).Value = 2
$ws.Cells.Item(402, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(402, 3).Value = "Coquimbo"
$ws.Cells.Item(402, 4).Value = 44783
$ws.Cells.Item(402, 4).Style = $ws.Cells.Item(401, 4).Style
$ws.Cells.Item(402, 4).NumberFormat = $ws.Cells.Item(401, 4).NumberFormat
$ws.Cells.Item(402, 5).Value = 4
$ws.Cells.Item(402, 6).Value = "Fruta"
$ws.Cells.Item(402, 7).Value = 100101
$ws.Cells.Item(402, 8).Value = "Berries"
$ws.Cells.Item(402, 9).Value = 100112025
$ws.Cells.Item(402, 10).Value = "Frutilla"
$ws.Cells.Item(402, 11).Value = "Sin especificar"
$ws.Cells.Item(402, 12).Value = "Primera"
$ws.Cells.Item(402, 13).Value = 160
$ws.Cells.Item(402, 14).Value = 25000
$ws.Cells.Item(402, 15).Value = 26000
$ws.Cells.Item(402, 16).Value = 25500
$ws.Cells.Item(402, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(402, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(402, 19).Value = 3643
$ws.Cells.Item(402, 20).Value = 7

# Row 403: new "Segunda" record dated 2022-08-10 (serial 44783)
$ws.Cells.Item(403, 1).Value = 2
$ws.Cells.Item(403, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(403, 3).Value = "Coquimbo"
$ws.Cells.Item(403, 4).Value = 44783
$ws.Cells.Item(403, 4).Style = $ws.Cells.Item(401, 4).Style
$ws.Cells.Item(403, 4).NumberFormat = $ws.Cells.Item(401, 4).NumberFormat
$ws.Cells.Item(403, 5).Value = 4
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100101
$ws.Cells.Item(403, 8).Value = "Berries"
$ws.Cells.Item(403, 9).Value = 100112025
$ws.Cells.Item(403, 10).Value = "Frutilla"
$ws.Cells.Item(403, 11).Value = "Sin especificar"
$ws.Cells.Item(403, 12).Value = "Segunda"
$ws.Cells.Item(403, 13).Value = 100
$ws.Cells.Item(403, 14).Value = 21000
$ws.Cells.Item(403, 15).Value = 22000
$ws.Cells.Item(403, 16).Value = 21500
$ws.Cells.Item(403, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(403, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(403, 19).Value = 3071
$ws.Cells.Item(403, 20).Value = 7
